$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 169.99634541439548
$ws.Range("C2").Value = 138.60883575451572
$ws.Range("D2").Value = 170.82698739089955
$ws.Range("E2").Value = 135.53889455340507

$ws.Range("B3").Value = 153.76991952659552
$ws.Range("C3").Value = 129.21433872432971
$ws.Range("D3").Value = 159.27475562953444
$ws.Range("E3").Value = 132.60611173301078

$ws.Range("B1:E3").Select()
